$d = $word.ActiveDocument
$table = $d.Tables(1)

function Set-CellLattice($table, $row, $col, $prob, $header, $r1, $r2) {
  $cell = $table.Cell($row, $col)
  $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + `
    '<w:t>' + $prob + '</w:t><w:br/>' + `
    '<w:t xml:space="preserve">' + $header + '</w:t><w:br/>' + `
    '<w:t xml:space="preserve">  ----</w:t><w:br/>' + `
    '<w:t>' + $r1 + '</w:t><w:br/>' + `
    '<w:t>' + $r2 + '</w:t>' + `
    '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  [void]$cell.Range.InsertXML($xml)
}

Set-CellLattice $table 1 1 "33 x 32" "  3    2" "3|    |" "3|    |"
Set-CellLattice $table 1 2 "21 x 81" "  8    1" "2|    |" "1|    |"
Set-CellLattice $table 1 3 "12 x 13" "  1    3" "1|    |" "2|    |"
Set-CellLattice $table 2 1 "18 x 27" "  2    7" "1|    |" "8|    |"
Set-CellLattice $table 2 2 "24 x 95" "  9    5" "2|    |" "4|    |"
Set-CellLattice $table 2 3 "74 x 37" "  3    7" "7|    |" "4|    |"
Set-CellLattice $table 3 1 "40 x 54" "  5    4" "4|    |" "0|    |"
Set-CellLattice $table 3 2 "85 x 39" "  3    9" "8|    |" "5|    |"
Set-CellLattice $table 3 3 "32 x 73" "  7    3" "3|    |" "2|    |"
Set-CellLattice $table 4 1 "99 x 85" "  8    5" "9|    |" "9|    |"
Set-CellLattice $table 4 2 "63 x 33" "  3    3" "6|    |" "3|    |"
Set-CellLattice $table 4 3 "28 x 66" "  6    6" "2|    |" "8|    |"
Set-CellLattice $table 5 1 "54 x 91" "  9    1" "5|    |" "4|    |"
Set-CellLattice $table 5 2 "23 x 84" "  8    4" "2|    |" "3|    |"
Set-CellLattice $table 5 3 "92 x 95" "  9    5" "9|    |" "2|    |"

Write-Output "done"
